{"js": "// Apply the same text substitutions as the target diff:\n//  - header date line\n//  - each \"NNN\u00d7N=\" problem cell, in document order\n// A couple of source strings (\"689\u00d74=\") occur twice in the document\n// with different replacement targets, so replacements are applied one at a\n// time (each immediately followed by context.sync()) and always target the\n// first remaining match - see the loop comment below for why that's safe.\n\nconst replacements = [\n  { find: \"2024-03-01 Friday\", replace: \"2024-03-02 Saturday\" },\n  { find: \"580\u00d77=\", replace: \"774\u00d79=\" },\n  { find: \"691\u00d78=\", replace: \"518\u00d73=\" },\n  { find: \"767\u00d74=\", replace: \"688\u00d73=\" },\n  { find: \"938\u00d73=\", replace: \"588\u00d74=\" },\n  { find: \"749\u00d78=\", replace: \"810\u00d75=\" },\n  { find: \"689\u00d74=\", replace: \"932\u00d73=\" }, // 1st occurrence\n  { find: \"264\u00d75=\", replace: \"334\u00d74=\" },\n  { find: \"676\u00d77=\", replace: \"662\u00d73=\" },\n  { find: \"555\u00d77=\", replace: \"542\u00d73=\" },\n  { find: \"584\u00d77=\", replace: \"150\u00d72=\" },\n  { find: \"929\u00d74=\", replace: \"437\u00d76=\" },\n  { find: \"854\u00d79=\", replace: \"799\u00d74=\" },\n  { find: \"893\u00d74=\", replace: \"854\u00d74=\" },\n  { find: \"913\u00d74=\", replace: \"819\u00d79=\" },\n  { find: \"490\u00d78=\", replace: \"317\u00d75=\" },\n  { find: \"147\u00d79=\", replace: \"917\u00d72=\" },\n  { find: \"430\u00d72=\", replace: \"521\u00d77=\" },\n  { find: \"689\u00d74=\", replace: \"825\u00d75=\" }, // 2nd occurrence\n  { find: \"328\u00d77=\", replace: \"792\u00d79=\" },\n  { find: \"508\u00d73=\", replace: \"745\u00d78=\" },\n  { find: \"941\u00d78=\", replace: \"644\u00d75=\" },\n  { find: \"286\u00d75=\", replace: \"659\u00d77=\" },\n  { find: \"588\u00d72=\", replace: \"964\u00d72=\" },\n  { find: \"695\u00d77=\", replace: \"437\u00d77=\" },\n  { find: \"619\u00d73=\", replace: \"984\u00d78=\" },\n];\n\n// Each replacement is applied (and synced) immediately, one at a time, in\n// document order. Because an already-applied replacement removes that\n// occurrence of the search text, the *next* remaining match for a repeated\n// search string (e.g. \"689\u00d74=\") is always back at index 0 - so we always\n// take the first hit rather than trying to pre-compute an occurrence index.\nfor (const { find, replace } of replacements) {\n  const results = context.document.body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length > 0) {\n    results.items[0].insertText(replace, \"Replace\");\n    await context.sync();\n  }\n}\n", "ps1": "# Apply the same text substitutions as the target diff:\n#  - header date line\n#  - each \"NNN\u00d7N=\" arithmetic-problem cell, in document order\n# A couple of source strings (\"689\u00d74=\") occur twice in the document with\n# different replacement targets, so each pair is applied as a single\n# \"replace just the next occurrence\" (wdReplaceOne) operation, in the same\n# order the values appear in the document - that way the first remaining\n# match is always the correct one, even for a repeated search string.\n\n$wdReplaceOne = 1\n$wdFindContinue = 1\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2024-03-01 Friday\", \"2024-03-02 Saturday\"),\n  @(\"580\u00d77=\", \"774\u00d79=\"),\n  @(\"691\u00d78=\", \"518\u00d73=\"),\n  @(\"767\u00d74=\", \"688\u00d73=\"),\n  @(\"938\u00d73=\", \"588\u00d74=\"),\n  @(\"749\u00d78=\", \"810\u00d75=\"),\n  @(\"689\u00d74=\", \"932\u00d73=\"),\n  @(\"264\u00d75=\", \"334\u00d74=\"),\n  @(\"676\u00d77=\", \"662\u00d73=\"),\n  @(\"555\u00d77=\", \"542\u00d73=\"),\n  @(\"584\u00d77=\", \"150\u00d72=\"),\n  @(\"929\u00d74=\", \"437\u00d76=\"),\n  @(\"854\u00d79=\", \"799\u00d74=\"),\n  @(\"893\u00d74=\", \"854\u00d74=\"),\n  @(\"913\u00d74=\", \"819\u00d79=\"),\n  @(\"490\u00d78=\", \"317\u00d75=\"),\n  @(\"147\u00d79=\", \"917\u00d72=\"),\n  @(\"430\u00d72=\", \"521\u00d77=\"),\n  @(\"689\u00d74=\", \"825\u00d75=\"),\n  @(\"328\u00d77=\", \"792\u00d79=\"),\n  @(\"508\u00d73=\", \"745\u00d78=\"),\n  @(\"941\u00d78=\", \"644\u00d75=\"),\n  @(\"286\u00d75=\", \"659\u00d77=\"),\n  @(\"588\u00d72=\", \"964\u00d72=\"),\n  @(\"695\u00d77=\", \"437\u00d77=\"),\n  @(\"619\u00d73=\", \"984\u00d78=\")\n)\n\nforeach ($pair in $pairs) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Execute($findText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $replaceText, $wdReplaceOne)\n}\n"}
